$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 "Marking": Right marks 5 -> 4, Wrong marks -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": total score 95 -> 76, and the "X / Y" summary text
$ws.Range("B12").Value = 76
$ws.Range("E12").Value = "76 / 112"
